$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 39628.5  # H3
$ws.Cells.Item(3, 10).Value = 39628.5  # J3
$ws.Cells.Item(3, 12).Value = 39628.5  # L3
$ws.Cells.Item(3, 14).Value = -39856.5  # N3
$ws.Cells.Item(33, 8).Value = 239.09375  # H33
$ws.Cells.Item(33, 9).Value = 143.04  # I33
$ws.Cells.Item(33, 10).Value = 582.1429000000001  # J33
$ws.Cells.Item(33, 11).Value = 143.04  # K33
$ws.Cells.Item(33, 12).Value = 582.1429000000001  # L33
$ws.Cells.Item(33, 13).Value = 85.96000000000001  # M33
$ws.Cells.Item(33, 14).Value = -1040.1429  # N33
$ws.Cells.Item(62, 8).Value = 3635  # H62
$ws.Cells.Item(62, 9).Value = 3635  # I62
$ws.Cells.Item(62, 10).Value = 0  # J62
$ws.Cells.Item(62, 11).Value = 3635  # K62
$ws.Cells.Item(62, 12).Value = 0  # L62
$ws.Cells.Item(62, 13).Value = -3011  # M62
$ws.Cells.Item(62, 14).ClearContents()  # N62
$ws.Cells.Item(65, 8).Value = 3635  # H65
$ws.Cells.Item(65, 9).Value = 3635  # I65
$ws.Cells.Item(65, 10).Value = 0  # J65
$ws.Cells.Item(65, 11).Value = 18175  # K65
$ws.Cells.Item(65, 12).Value = 0  # L65
$ws.Cells.Item(65, 13).Value = -15055  # M65
$ws.Cells.Item(65, 14).Value = -15055  # N65
$ws.Cells.Item(102, 8).Value = 39628.5  # H102
$ws.Cells.Item(102, 10).Value = 39628.5  # J102
$ws.Cells.Item(102, 12).Value = 39628.5  # L102
$ws.Cells.Item(102, 14).Value = -46118.5  # N102
$ws.Cells.Item(111, 8).Value = 2438.16  # H111
$ws.Cells.Item(111, 9).Value = 2001.2941  # I111
$ws.Cells.Item(111, 10).Value = 3366.5  # J111
$ws.Cells.Item(111, 11).Value = 6003.8823  # K111
$ws.Cells.Item(111, 12).Value = 10099.5  # L111
$ws.Cells.Item(111, 13).Value = -2936.8823  # M111
$ws.Cells.Item(111, 14).Value = -16233.5  # N111

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4880.4443  # H61
$ws.Cells.Item(61, 9).Value = 3449.0571  # I61
$ws.Cells.Item(61, 10).Value = 9890.299999999999  # J61
$ws.Cells.Item(61, 11).Value = 3449.0571  # K61
$ws.Cells.Item(61, 12).Value = 9890.299999999999  # L61
$ws.Cells.Item(61, 13).Value = -3237.0571  # M61
$ws.Cells.Item(61, 14).Value = -10314.3  # N61
$ws.Cells.Item(74, 8).Value = 5379.241  # H74
$ws.Cells.Item(74, 9).Value = 6062.8  # I74
$ws.Cells.Item(74, 10).Value = 1107  # J74
$ws.Cells.Item(74, 11).Value = 6062.8  # K74
$ws.Cells.Item(74, 12).Value = 1107  # L74
$ws.Cells.Item(74, 13).Value = -5188.8  # M74
$ws.Cells.Item(74, 14).Value = -2855  # N74
$ws.Cells.Item(77, 8).Value = 5379.241  # H77
$ws.Cells.Item(77, 9).Value = 6062.8  # I77
$ws.Cells.Item(77, 10).Value = 1107  # J77
$ws.Cells.Item(77, 11).Value = 30314  # K77
$ws.Cells.Item(77, 12).Value = 5535  # L77
$ws.Cells.Item(77, 13).Value = -25946  # M77
$ws.Cells.Item(77, 14).Value = -14271  # N77
$ws.Cells.Item(122, 8).Value = 1229.5652  # H122
$ws.Cells.Item(122, 9).Value = 1030.625  # I122
$ws.Cells.Item(122, 10).Value = 1684.2858  # J122
$ws.Cells.Item(122, 11).Value = 3091.875  # K122
$ws.Cells.Item(122, 12).Value = 5052.857400000001  # L122
$ws.Cells.Item(122, 13).Value = -641.875  # M122
$ws.Cells.Item(122, 14).Value = -9952.857400000001  # N122
$ws.Cells.Item(132, 8).Value = 1945.1459  # H132
$ws.Cells.Item(132, 9).Value = 1076.9333  # I132
$ws.Cells.Item(132, 10).Value = 3392.1667  # J132
$ws.Cells.Item(132, 11).Value = 3230.7999  # K132
$ws.Cells.Item(132, 12).Value = 10176.5001  # L132
$ws.Cells.Item(132, 13).Value = -700.7999  # M132
$ws.Cells.Item(132, 14).Value = -15236.5001  # N132
$ws.Cells.Item(136, 8).Value = 4880.4443  # H136
$ws.Cells.Item(136, 9).Value = 3449.0571  # I136
$ws.Cells.Item(136, 10).Value = 9890.299999999999  # J136
$ws.Cells.Item(136, 11).Value = 10347.1713  # K136
$ws.Cells.Item(136, 12).Value = 29670.9  # L136
$ws.Cells.Item(136, 13).Value = -7797.1713  # M136
$ws.Cells.Item(136, 14).Value = -34770.89999999999  # N136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 18548.4  # H26
$ws.Cells.Item(26, 9).Value = 15685.5  # I26
$ws.Cells.Item(26, 10).Value = 30000  # J26
$ws.Cells.Item(26, 11).Value = 15685.5  # K26
$ws.Cells.Item(26, 12).Value = 30000  # L26
$ws.Cells.Item(26, 13).Value = -15393.5  # M26
$ws.Cells.Item(26, 14).Value = -30584  # N26
$ws.Cells.Item(57, 8).Value = 50886.668  # H57
$ws.Cells.Item(57, 10).Value = 50886.668  # J57
$ws.Cells.Item(57, 12).Value = 50886.668  # L57
$ws.Cells.Item(57, 14).Value = -52326.668  # N57
$ws.Cells.Item(107, 8).Value = 1320.25  # H107
$ws.Cells.Item(107, 9).Value = 938.5  # I107
$ws.Cells.Item(107, 10).Value = 1702  # J107
$ws.Cells.Item(107, 11).Value = 938.5  # K107
$ws.Cells.Item(107, 12).Value = 1702  # L107
$ws.Cells.Item(107, 13).Value = 981.5  # M107
$ws.Cells.Item(107, 14).Value = -5542  # N107
$ws.Cells.Item(134, 8).Value = 2162.9443  # H134
$ws.Cells.Item(134, 9).Value = 2038.6428  # I134
$ws.Cells.Item(134, 11).Value = 6115.928400000001  # K134
$ws.Cells.Item(134, 13).Value = -3580.928400000001  # M134
$ws.Cells.Item(136, 8).Value = 50886.668  # H136
$ws.Cells.Item(136, 10).Value = 50886.668  # J136
$ws.Cells.Item(136, 12).Value = 50886.668  # L136
$ws.Cells.Item(136, 14).Value = -61086.668  # N136

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 0  # H28
$ws.Cells.Item(28, 10).Value = 0  # J28
$ws.Cells.Item(28, 12).Value = 0  # L28
$ws.Cells.Item(28, 14).ClearContents()  # N28
$ws.Cells.Item(58, 8).Value = 2022247.4  # H58
$ws.Cells.Item(58, 9).Value = 5348756.5  # I58
$ws.Cells.Item(58, 10).Value = 2581.2144  # J58
$ws.Cells.Item(58, 11).Value = 5348756.5  # K58
$ws.Cells.Item(58, 12).Value = 2581.2144  # L58
$ws.Cells.Item(58, 13).Value = -5348553.5  # M58
$ws.Cells.Item(58, 14).Value = -2987.2144  # N58
$ws.Cells.Item(107, 8).Value = 700.2895  # H107
$ws.Cells.Item(107, 9).Value = 706.7619  # I107
$ws.Cells.Item(107, 10).Value = 692.2941  # J107
$ws.Cells.Item(107, 11).Value = 706.7619  # K107
$ws.Cells.Item(107, 12).Value = 692.2941  # L107
$ws.Cells.Item(107, 13).Value = 1213.2381  # M107
$ws.Cells.Item(107, 14).Value = -4532.2941  # N107
$ws.Cells.Item(122, 8).Value = 14382.363  # H122
$ws.Cells.Item(122, 9).Value = 6899  # I122
$ws.Cells.Item(122, 10).Value = 34338  # J122
$ws.Cells.Item(122, 11).Value = 20697  # K122
$ws.Cells.Item(122, 12).Value = 103014  # L122
$ws.Cells.Item(122, 13).Value = -18247  # M122
$ws.Cells.Item(122, 14).Value = -107914  # N122
$ws.Cells.Item(132, 8).Value = 3974.3542  # H132
$ws.Cells.Item(132, 9).Value = 4596.5483  # I132
$ws.Cells.Item(132, 10).Value = 2839.7646  # J132
$ws.Cells.Item(132, 11).Value = 13789.6449  # K132
$ws.Cells.Item(132, 12).Value = 8519.293799999999  # L132
$ws.Cells.Item(132, 13).Value = -11259.6449  # M132
$ws.Cells.Item(132, 14).Value = -13579.2938  # N132
$ws.Cells.Item(134, 8).Value = 1900.3615  # H134
$ws.Cells.Item(134, 9).Value = 1044.0714  # I134
$ws.Cells.Item(134, 11).Value = 3132.2142  # K134
$ws.Cells.Item(134, 13).Value = -597.2142000000003  # M134
$ws.Cells.Item(136, 8).Value = 2022247.4  # H136
$ws.Cells.Item(136, 9).Value = 5348756.5  # I136
$ws.Cells.Item(136, 10).Value = 2581.2144  # J136
$ws.Cells.Item(136, 11).Value = 16046269.5  # K136
$ws.Cells.Item(136, 12).Value = 7743.6432  # L136
$ws.Cells.Item(136, 13).Value = -16043719.5  # M136
$ws.Cells.Item(136, 14).Value = -12843.6432  # N136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 16677154  # H5
$ws.Cells.Item(5, 9).Value = 608.125  # I5
$ws.Cells.Item(5, 10).Value = 83383336  # J5
$ws.Cells.Item(5, 11).Value = 1824.375  # K5
$ws.Cells.Item(5, 12).Value = 250150008  # L5
$ws.Cells.Item(5, 13).Value = -1712.375  # M5
$ws.Cells.Item(5, 14).Value = -250150232  # N5
$ws.Cells.Item(22, 8).Value = 1587.7778  # H22
$ws.Cells.Item(22, 9).Value = 683.3333  # I22
$ws.Cells.Item(22, 10).Value = 2040  # J22
$ws.Cells.Item(22, 11).Value = 2049.9999  # K22
$ws.Cells.Item(22, 12).Value = 6120  # L22
$ws.Cells.Item(22, 13).Value = -1880.9999  # M22
$ws.Cells.Item(22, 14).Value = -6458  # N22
$ws.Cells.Item(23, 8).Value = 818.7143  # H23
$ws.Cells.Item(23, 10).Value = 125.454544  # J23
$ws.Cells.Item(23, 12).Value = 376.363632  # L23
$ws.Cells.Item(23, 14).Value = -846.3636320000001  # N23
$ws.Cells.Item(27, 8).Value = 1587.7778  # H27
$ws.Cells.Item(27, 9).Value = 683.3333  # I27
$ws.Cells.Item(27, 10).Value = 2040  # J27
$ws.Cells.Item(27, 11).Value = 2049.9999  # K27
$ws.Cells.Item(27, 12).Value = 6120  # L27
$ws.Cells.Item(27, 13).Value = -1947.9999  # M27
$ws.Cells.Item(27, 14).Value = -6324  # N27
$ws.Cells.Item(33, 8).Value = 100.625  # H33
$ws.Cells.Item(33, 10).Value = 111.42857  # J33
$ws.Cells.Item(33, 12).Value = 668.57142  # L33
$ws.Cells.Item(33, 14).Value = -1234.57142  # N33
$ws.Cells.Item(40, 8).Value = 140.83333  # H40
$ws.Cells.Item(40, 9).Value = 41.666668  # I40
$ws.Cells.Item(40, 10).Value = 240  # J40
$ws.Cells.Item(40, 11).Value = 166.666672  # K40
$ws.Cells.Item(40, 12).Value = 960  # L40
$ws.Cells.Item(40, 13).Value = -97.66667200000001  # M40
$ws.Cells.Item(40, 14).Value = -1098  # N40
$ws.Cells.Item(64, 8).Value = 83335736  # H64
$ws.Cells.Item(64, 9).Value = 250000750  # I64
$ws.Cells.Item(64, 10).Value = 3224.375  # J64
$ws.Cells.Item(64, 11).Value = 750002250  # K64
$ws.Cells.Item(64, 12).Value = 9673.125  # L64
$ws.Cells.Item(64, 13).Value = -750001980  # M64
$ws.Cells.Item(64, 14).Value = -10213.125  # N64
$ws.Cells.Item(67, 8).Value = 83335736  # H67
$ws.Cells.Item(67, 9).Value = 250000750  # I67
$ws.Cells.Item(67, 10).Value = 3224.375  # J67
$ws.Cells.Item(67, 11).Value = 750002250  # K67
$ws.Cells.Item(67, 12).Value = 9673.125  # L67
$ws.Cells.Item(67, 13).Value = -750001314  # M67
$ws.Cells.Item(67, 14).Value = -11545.125  # N67
$ws.Cells.Item(122, 8).Value = 851.3333  # H122
$ws.Cells.Item(122, 10).Value = 990.1539  # J122
$ws.Cells.Item(122, 12).Value = 8911.3851  # L122
$ws.Cells.Item(122, 14).Value = -13811.3851  # N122
$ws.Cells.Item(135, 8).Value = 16677154  # H135
$ws.Cells.Item(135, 9).Value = 608.125  # I135
$ws.Cells.Item(135, 10).Value = 83383336  # J135
$ws.Cells.Item(135, 11).Value = 5473.125  # K135
$ws.Cells.Item(135, 12).Value = 750450024  # L135
$ws.Cells.Item(135, 13).Value = -2938.125  # M135
$ws.Cells.Item(135, 14).Value = -750455094  # N135

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 464.64285  # H107
$ws.Cells.Item(107, 9).Value = 181.25  # I107
$ws.Cells.Item(107, 11).Value = 181.25  # K107
$ws.Cells.Item(107, 13).Value = 1738.75  # M107
$ws.Cells.Item(122, 8).Value = 3311.1924  # H122
$ws.Cells.Item(122, 9).Value = 4488  # I122
$ws.Cells.Item(122, 10).Value = 1706.4546  # J122
$ws.Cells.Item(122, 11).Value = 13464  # K122
$ws.Cells.Item(122, 12).Value = 5119.3638  # L122
$ws.Cells.Item(122, 13).Value = -11014  # M122
$ws.Cells.Item(122, 14).Value = -10019.3638  # N122
$ws.Cells.Item(132, 8).Value = 2274.1667  # H132
$ws.Cells.Item(132, 9).Value = 2153.348  # I132
$ws.Cells.Item(132, 10).Value = 2671.1428  # J132
$ws.Cells.Item(132, 11).Value = 6460.044  # K132
$ws.Cells.Item(132, 12).Value = 8013.428400000001  # L132
$ws.Cells.Item(132, 13).Value = -3930.044  # M132
$ws.Cells.Item(132, 14).Value = -13073.4284  # N132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(81, 8).Value = 40181  # H81
$ws.Cells.Item(81, 10).Value = 40181  # J81
$ws.Cells.Item(81, 12).Value = 40181  # L81
$ws.Cells.Item(81, 14).Value = -42177  # N81
$ws.Cells.Item(84, 8).Value = 40181  # H84
$ws.Cells.Item(84, 10).Value = 40181  # J84
$ws.Cells.Item(84, 12).Value = 120543  # L84
$ws.Cells.Item(84, 14).Value = -130527  # N84
$ws.Cells.Item(132, 8).Value = 5454  # H132
$ws.Cells.Item(132, 9).Value = 7308.6  # I132
$ws.Cells.Item(132, 10).Value = 2981.2  # J132
$ws.Cells.Item(132, 11).Value = 21925.8  # K132
$ws.Cells.Item(132, 12).Value = 8943.599999999999  # L132
$ws.Cells.Item(132, 13).Value = -19395.8  # M132
$ws.Cells.Item(132, 14).Value = -14003.6  # N132
$ws.Cells.Item(136, 8).Value = 4091.3333  # H136
$ws.Cells.Item(136, 9).Value = 2101.3872  # I136
$ws.Cells.Item(136, 10).Value = 7175.75  # J136
$ws.Cells.Item(136, 11).Value = 6304.1616  # K136
$ws.Cells.Item(136, 12).Value = 21527.25  # L136
$ws.Cells.Item(136, 13).Value = -3754.1616  # M136
$ws.Cells.Item(136, 14).Value = -26627.25  # N136
